$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 stays the same (jEgFPF3LKv@gmail.com with mailto hyperlink)

# Re-assert A3 (same text as before) and add new row A4
$ws.Range("A3").Value = "MgXMUk5Brz@gmail.com"
$ws.Range("A4").Value = "VsPxaWd40T@gmail.com"

# Update header cell (A1)
$ws.Range("A1").Value = "Registered_EmailIds"

# Append remaining new rows
$ws.Range("A5").Value = "MiCcv4AjCP@gmail.com"
$ws.Range("A6").Value = "SgQCS3JTB6@gmail.com"

# Move selection back to A1 (removes stale B2 selection reference)
$ws.Range("A1").Select()
